$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.849.29"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.649.20"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "1.642.91"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "26.849.65"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "0.0₃0737"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("E21").Value = "  +11.19%  "
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0510"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "1.297.75"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("D43").Value = "1.784.61"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0521"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0974"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
